$wb = $excel.ActiveWorkbook

# Target stored column width (OOXML <col width=.../>) is 17.2159881591797.
# The COM ColumnWidth setter snaps onto this engine's character-width pixel
# grid (steps of 1/6), so 16.33 is the input that lands on the closest
# reachable stored width (17.1666..., the nearest grid point to the target).
$narrowColumnWidth = 16.33

# --- Sheet "Overview": status text + handoff timestamp ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-24 09:02:41"
$ws1.Columns.Item(5).ColumnWidth = $narrowColumnWidth
$ws1.Columns.Item(6).ColumnWidth = $narrowColumnWidth

# --- Sheet "zh-cn": status text + handoff timestamp ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-08-24 09:02:36"
$ws2.Columns.Item(3).ColumnWidth = $narrowColumnWidth

# --- Sheet "de-de": status text + handoff timestamp ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-08-24 09:02:41"
$ws3.Columns.Item(3).ColumnWidth = $narrowColumnWidth
